# Applies the "Ripple" template update:
#  - Adds a new "Assay" worksheet with Setting/Value rows
#  - Tweaks sheetView selections on Patterns / Compounds sheets
#  - Removes the (now unused) duplicate cell style from row 1 of Compounds

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Assay" worksheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

$assayData = @(
    @("Setting", "Value"),
    @("DMSO Tolerance", 0.005),
    @("Well Volume (µL)", 25),
    @("Backfill (µL)", 10),
    @("Allowed Error", 0.1),
    @("Destination Replicates", 1),
    @("Use Intermediate Plates", 1),
    @("DMSO Normalization", 1)
)

for ($i = 0; $i -lt $assayData.Count; $i++) {
    $rowNum = $i + 1
    $assay.Cells.Item($rowNum, 1).Value = $assayData[$i][0]
    $assay.Cells.Item($rowNum, 2).Value = $assayData[$i][1]
}

[void]$assay.Range("A1:B8").Select()

# ---------------------------------------------------------------------------
# 2. Patterns sheet: becomes the active tab, selection moves to D12
# ---------------------------------------------------------------------------
$patterns = $wb.Worksheets.Item("Patterns")
[void]$patterns.Select()
[void]$patterns.Range("D12").Select()

# ---------------------------------------------------------------------------
# 3. Compounds sheet: remove leftover duplicate style from header row,
#    selection moves back to its previous cell (M26), no longer the active tab
# ---------------------------------------------------------------------------
$compounds = $wb.Worksheets.Item("Compounds")
$compounds.Range("A1:F1").Style = "Normal"
[void]$compounds.Select()
[void]$compounds.Range("M26").Select()

# Re-activate Patterns last so it ends up as the active/visible tab
[void]$patterns.Select()
